$p = $ppt.ActivePresentation

# Slide 1: "Header" + " " + "with" + " " + "inline code" -> "Header with " + "inline code"
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$tr1.Characters(1, 12).Text = "Header with "

# Slide 2: "Syntax" + " " + "highlighting" -> "Syntax highlighting"
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$tr2.Characters(1, 20).Text = "Syntax highlighting"

# Slide 3: "Two" + " " + "column" + " " + "slide" -> "Two column slide"
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tr3.Characters(1, 17).Text = "Two column slide"
